$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy date-cell formatting (style index 2 / border+bold+center+numFmt) from the last
# existing row (A385) down through the new rows so the new date cells in column A
# keep the same look as all the other date cells above them.
$ws.Range("A385").Copy() | Out-Null
$ws.Range("A386:A464").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Populate the new rows (data update "fino a 8/12" / through 2021-12-08)
$data = @(
    ,@(44460, 2, 20, 113.836871762764)
    ,@(44461, 1, 21, 119.5287153509022)
    ,@(44462, 1, 15, 85.37765382207297)
    ,@(44463, 3, 16, 91.06949741021117)
    ,@(44464, 1, 14, 79.68581023393477)
    ,@(44465, 0, 10, 56.91843588138198)
    ,@(44466, 3, 11, 62.61027946952018)
    ,@(44467, 1, 10, 56.91843588138198)
    ,@(44468, 0, 9, 51.22659229324378)
    ,@(44469, 0, 8, 45.53474870510559)
    ,@(44470, 0, 5, 28.45921794069099)
    ,@(44471, 2, 6, 34.15106152882919)
    ,@(44472, 0, 6, 34.15106152882919)
    ,@(44473, 2, 5, 28.45921794069099)
    ,@(44474, 2, 6, 34.15106152882919)
    ,@(44475, 0, 6, 34.15106152882919)
    ,@(44476, 1, 7, 39.84290511696739)
    ,@(44477, 0, 7, 39.84290511696739)
    ,@(44478, 0, 5, 28.45921794069099)
    ,@(44479, 1, 6, 34.15106152882919)
    ,@(44480, 1, 5, 28.45921794069099)
    ,@(44481, 0, 3, 17.07553076441459)
    ,@(44482, 1, 4, 22.76737435255279)
    ,@(44483, 0, 3, 17.07553076441459)
    ,@(44484, 0, 3, 17.07553076441459)
    ,@(44485, 0, 3, 17.07553076441459)
    ,@(44486, 0, 2, 11.3836871762764)
    ,@(44487, 1, 2, 11.3836871762764)
    ,@(44488, 1, 3, 17.07553076441459)
    ,@(44489, 1, 3, 17.07553076441459)
    ,@(44490, 2, 5, 28.45921794069099)
    ,@(44491, 3, 8, 45.53474870510559)
    ,@(44492, 0, 8, 45.53474870510559)
    ,@(44493, 1, 9, 51.22659229324378)
    ,@(44494, 0, 8, 45.53474870510559)
    ,@(44495, 0, 7, 39.84290511696739)
    ,@(44496, 0, 6, 34.15106152882919)
    ,@(44497, 0, 4, 22.76737435255279)
    ,@(44498, 1, 2, 11.3836871762764)
    ,@(44499, 1, 3, 17.07553076441459)
    ,@(44500, 0, 2, 11.3836871762764)
    ,@(44501, 4, 6, 34.15106152882919)
    ,@(44502, 0, 6, 34.15106152882919)
    ,@(44503, 0, 6, 34.15106152882919)
    ,@(44504, 2, 8, 45.53474870510559)
    ,@(44505, 1, 8, 45.53474870510559)
    ,@(44506, 3, 10, 56.91843588138198)
    ,@(44507, 0, 10, 56.91843588138198)
    ,@(44508, 1, 7, 39.84290511696739)
    ,@(44509, 2, 9, 51.22659229324378)
    ,@(44510, 0, 9, 51.22659229324378)
    ,@(44511, 1, 8, 45.53474870510559)
    ,@(44512, 4, 11, 62.61027946952018)
    ,@(44513, 2, 10, 56.91843588138198)
    ,@(44514, 2, 12, 68.30212305765838)
    ,@(44515, 0, 11, 62.61027946952018)
    ,@(44516, 15, 24, 136.6042461153168)
    ,@(44517, 0, 24, 136.6042461153168)
    ,@(44518, 8, 31, 176.4471512322841)
    ,@(44519, 8, 35, 199.2145255848369)
    ,@(44520, 8, 41, 233.3655871136661)
    ,@(44521, 5, 44, 250.4411178780807)
    ,@(44522, 6, 50, 284.5921794069099)
    ,@(44523, 5, 40, 227.6737435255279)
    ,@(44524, 5, 45, 256.1329614662189)
    ,@(44525, 6, 43, 244.7492742899425)
    ,@(44526, 9, 44, 250.4411178780807)
    ,@(44527, 2, 38, 216.2900563492515)
    ,@(44528, 2, 35, 199.2145255848369)
    ,@(44529, 8, 37, 210.5982127611133)
    ,@(44530, 7, 39, 221.9818999373897)
    ,@(44531, 0, 34, 193.5226819966987)
    ,@(44532, 6, 34, 193.5226819966987)
    ,@(44533, 2, 27, 153.6797768797313)
    ,@(44534, 4, 29, 165.0634640560077)
    ,@(44535, 3, 30, 170.7553076441459)
    ,@(44536, 8, 30, 170.7553076441459)
    ,@(44537, 11, 34, 193.5226819966987)
    ,@(44538, 4, 38, 216.2900563492515)
)

$startRow = 386
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
